$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (id 17) - Capital city
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Capital city"
$ws.Range("C18").Value = "Allows a player to designate one of their kingdoms per plane, as a Capital City. Capital Cities allows you to manage multiple kingdoms from one central kingdom."
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 2
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 1

# Row 19 (id 18) - Markets and Economy
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Markets and Economy"
$ws.Range("C19").Value = "Allows players to create Market Place in thier kingdoms to allow for resources to be ransfered."
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 3
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 1

# Row 20 (id 19) - Moving resources
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Moving resources"
$ws.Range("C20").Value = "Allows you to request a total of +90,000 resources at max level for a kingdom. By default a kingdom with a market can request, from another kingdom that also has a market, resources in the total of 5,000 or 10,000 is the kingdom being requested has airships. In total this allows you to mvoe 100,000 resources of a single or of all types."
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 2
$ws.Range("G20").Value = 10000
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 1
